$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "Mean mass flow rate" values (column B) from the re-run simulation
$ws.Range("B2").Value = 0.056033548096320121
$ws.Range("B3").Value = 0.40312088999951418
$ws.Range("B4").Value = 0.083988009950046294
$ws.Range("B5").Value = 0.031200612844029342
$ws.Range("B6").Value = 0.22097121656698404
$ws.Range("B7").Value = 0.079829303711251395
$ws.Range("B8").Value = 0.03662998542337341
$ws.Range("B9").Value = 0.3175151946015104

# Column widths were also nudged slightly wider
$ws.Columns.Item(1).ColumnWidth = 25.584635416666668
$ws.Columns.Item(2).ColumnWidth = 18.584635416666668
$ws.Columns.Item(3).ColumnWidth = 3.91796875
